# Update the Entsoe consumption forecast sheet:
#  - Column A (Forecasted Consumption (MW)) gets refreshed forecast values
#  - Column B (Timestamp) shifts forward by 4 days (new forecast run date)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newConsumption = @(5330,5280,5230,5180,5150,5110,5080,5060,5060,5050,5050,5050,5050,5050,5040,5050,5060,5090,5130,5200,5270,5360,5460,5560,5660,5750,5820,5870,5900,5900,5880,5830,5770,5690,5600,5510,5420,5340,5270,5210,5160,5120,5090,5070,5050,5040,5040,5040,5050,5070,5090,5120,5140,5170,5200,5230,5260,5290,5330,5380,5440,5510,5600,5680,5780,5880,5960,6060,6150,6260,6370,6470,6560,6650,6730,6800,6860,6920,6980,7060,7110,7100,7070,7030,6950,6840,6680,6500,6330,6160,6000,5860,5590,5550,5480,5410)

$dayShift = 4

for ($i = 0; $i -lt $newConsumption.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newConsumption[$i]
    $ws.Cells.Item($row, 2).Value = $ws.Cells.Item($row, 2).Value2 + $dayShift
}
